$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 18,13

$arr[0,0] = 187.0
$arr[0,1] = 'Monday, Jan 09'
$arr[0,2] = '2:30 PM'
$arr[0,3] = 'FR6252'
$arr[0,4] = 'Stockholm'
$arr[0,5] = '(ARN)'
$arr[0,6] = 'Ryanair '
$arr[0,7] = 'B38M'
$arr[0,8] = '(SP-RZL)'
$arr[0,9] = '2:44 PM'
$arr[0,10] = $null
$arr[0,11] = '0 hours, 14 minutes'
$arr[0,12] = $null

$arr[1,0] = 188.0
$arr[1,1] = 'Monday, Jan 09'
$arr[1,2] = '2:30 PM'
$arr[1,3] = 'LX1371'
$arr[1,4] = 'Zurich'
$arr[1,5] = '(ZRH)'
$arr[1,6] = 'Helvetic Airways '
$arr[1,7] = 'E295'
$arr[1,8] = '(HB-AZI)'
$arr[1,9] = '2:49 PM'
$arr[1,10] = $null
$arr[1,11] = '0 hours, 19 minutes'
$arr[1,12] = $null

$arr[2,0] = 189.0
$arr[2,1] = 'Monday, Jan 09'
$arr[2,2] = '2:30 PM'
$arr[2,3] = 'OS598'
$arr[2,4] = 'Vienna'
$arr[2,5] = '(VIE)'
$arr[2,6] = 'Austrian Airlines '
$arr[2,7] = 'E195'
$arr[2,8] = '(OE-LWQ)'
$arr[2,9] = '2:42 PM'
$arr[2,10] = $null
$arr[2,11] = '0 hours, 12 minutes'
$arr[2,12] = $null

$arr[3,0] = 190.0
$arr[3,1] = 'Monday, Jan 09'
$arr[3,2] = '2:45 PM'
$arr[3,3] = 'FR6877'
$arr[3,4] = 'Milan'
$arr[3,5] = '(BGY)'
$arr[3,6] = 'Buzz '
$arr[3,7] = 'B38M'
$arr[3,8] = '(SP-RZH)'
$arr[3,9] = '2:59 PM'
$arr[3,10] = $null
$arr[3,11] = '0 hours, 14 minutes'
$arr[3,12] = $null

$arr[4,0] = 191.0
$arr[4,1] = 'Monday, Jan 09'
$arr[4,2] = '3:00 PM'
$arr[4,3] = 'LO3906'
$arr[4,4] = 'Warsaw'
$arr[4,5] = '(WAW)'
$arr[4,6] = 'LOT '
$arr[4,7] = 'E190'
$arr[4,8] = '(SP-LMH)'
$arr[4,9] = '3:34 PM'
$arr[4,10] = $null
$arr[4,11] = '0 hours, 34 minutes'
$arr[4,12] = $null

$arr[5,0] = 192.0
$arr[5,1] = 'Monday, Jan 09'
$arr[5,2] = '3:00 PM'
$arr[5,3] = 'LPR41'
$arr[5,4] = 'Warsaw'
$arr[5,5] = '(WAW)'
$arr[5,6] = 'Polish Medical Air Rescue '
$arr[5,7] = 'LJ75'
$arr[5,8] = '(SP-MXR)'
$arr[5,9] = '5:04 PM'
$arr[5,10] = $null
$arr[5,11] = '2 hours, 4 minutes'
$arr[5,12] = $null

$arr[6,0] = 193.0
$arr[6,1] = 'Monday, Jan 09'
$arr[6,2] = '3:10 PM'
$arr[6,3] = 'FR7955'
$arr[6,4] = 'Prague'
$arr[6,5] = '(PRG)'
$arr[6,6] = 'Ryanair '
$arr[6,7] = 'B738'
$arr[6,8] = '(SP-RSY)'
$arr[6,9] = '3:32 PM'
$arr[6,10] = $null
$arr[6,11] = '0 hours, 22 minutes'
$arr[6,12] = $null

$arr[7,0] = 194.0
$arr[7,1] = 'Monday, Jan 09'
$arr[7,2] = '3:15 PM'
$arr[7,3] = 'FR6228'
$arr[7,4] = 'Tel Aviv'
$arr[7,5] = '(TLV)'
$arr[7,6] = 'Buzz '
$arr[7,7] = 'B38M'
$arr[7,8] = '(SP-RZD)'
$arr[7,9] = '3:30 PM'
$arr[7,10] = $null
$arr[7,11] = '0 hours, 15 minutes'
$arr[7,12] = $null

$arr[8,0] = 195.0
$arr[8,1] = 'Monday, Jan 09'
$arr[8,2] = '3:15 PM'
$arr[8,3] = 'FR6361'
$arr[8,4] = 'Shannon'
$arr[8,5] = '(SNN)'
$arr[8,6] = 'Ryanair '
$arr[8,7] = 'B738'
$arr[8,8] = '(EI-EBP)'
$arr[8,9] = '3:25 PM'
$arr[8,10] = $null
$arr[8,11] = '0 hours, 10 minutes'
$arr[8,12] = $null

$arr[9,0] = 196.0
$arr[9,1] = 'Monday, Jan 09'
$arr[9,2] = '3:40 PM'
$arr[9,3] = 'FR3722'
$arr[9,4] = 'Billund'
$arr[9,5] = '(BLL)'
$arr[9,6] = 'Ryanair '
$arr[9,7] = 'B738'
$arr[9,8] = '(9H-QBX)'
$arr[9,9] = '3:49 PM'
$arr[9,10] = $null
$arr[9,11] = '0 hours, 9 minutes'
$arr[9,12] = $null

$arr[10,0] = 197.0
$arr[10,1] = 'Monday, Jan 09'
$arr[10,2] = '3:40 PM'
$arr[10,3] = 'LG5742'
$arr[10,4] = 'Luxembourg'
$arr[10,5] = '(LUX)'
$arr[10,6] = 'Luxair '
$arr[10,7] = 'DH8D'
$arr[10,8] = '(LX-LGE)'
$arr[10,9] = '3:36 PM'
$arr[10,10] = $null
$arr[10,11] = '0 hours, -4 minutes'
$arr[10,12] = $null

$arr[11,0] = 198.0
$arr[11,1] = 'Monday, Jan 09'
$arr[11,2] = '3:50 PM'
$arr[11,3] = 'U23816'
$arr[11,4] = 'Paris'
$arr[11,5] = '(CDG)'
$arr[11,6] = 'easyJet '
$arr[11,7] = 'A320'
$arr[11,8] = '(OE-IVS)'
$arr[11,9] = '4:07 PM'
$arr[11,10] = $null
$arr[11,11] = '0 hours, 17 minutes'
$arr[11,12] = $null

$arr[12,0] = 199.0
$arr[12,1] = 'Monday, Jan 09'
$arr[12,2] = '4:10 PM'
$arr[12,3] = 'FR6248'
$arr[12,4] = 'Manchester'
$arr[12,5] = '(MAN)'
$arr[12,6] = 'Buzz '
$arr[12,7] = 'B38M'
$arr[12,8] = '(SP-RZA)'
$arr[12,9] = '4:23 PM'
$arr[12,10] = $null
$arr[12,11] = '0 hours, 13 minutes'
$arr[12,12] = $null

$arr[13,0] = 200.0
$arr[13,1] = 'Monday, Jan 09'
$arr[13,2] = '4:40 PM'
$arr[13,3] = 'FR2713'
$arr[13,4] = 'London'
$arr[13,5] = '(STN)'
$arr[13,6] = 'Ryanair '
$arr[13,7] = 'B738'
$arr[13,8] = '(EI-DWY)'
$arr[13,9] = '4:55 PM'
$arr[13,10] = $null
$arr[13,11] = '0 hours, 15 minutes'
$arr[13,12] = $null

$arr[14,0] = 201.0
$arr[14,1] = 'Monday, Jan 09'
$arr[14,2] = '4:40 PM'
$arr[14,3] = 'FR3054'
$arr[14,4] = 'Barcelona'
$arr[14,5] = '(BCN)'
$arr[14,6] = 'Buzz '
$arr[14,7] = 'B38M'
$arr[14,8] = '(SP-RZF)'
$arr[14,9] = '4:43 PM'
$arr[14,10] = $null
$arr[14,11] = '0 hours, 3 minutes'
$arr[14,12] = $null

$arr[15,0] = 202.0
$arr[15,1] = 'Monday, Jan 09'
$arr[15,2] = '4:50 PM'
$arr[15,3] = 'KL1996'
$arr[15,4] = 'Amsterdam'
$arr[15,5] = '(AMS)'
$arr[15,6] = 'KLM '
$arr[15,7] = 'E190'
$arr[15,8] = '(PH-EXC)'
$arr[15,9] = '5:08 PM'
$arr[15,10] = $null
$arr[15,11] = '0 hours, 18 minutes'
$arr[15,12] = $null

$arr[16,0] = 203.0
$arr[16,1] = 'Monday, Jan 09'
$arr[16,2] = '4:55 PM'
$arr[16,3] = 'W65051'
$arr[16,4] = 'Larnaca'
$arr[16,5] = '(LCA)'
$arr[16,6] = 'Wizz Air '
$arr[16,7] = 'A21N'
$arr[16,8] = '(HA-LZI)'
$arr[16,9] = '5:14 PM'
$arr[16,10] = $null
$arr[16,11] = '0 hours, 19 minutes'
$arr[16,12] = $null

$arr[17,0] = 204.0
$arr[17,1] = 'Monday, Jan 09'
$arr[17,2] = '5:10 PM'
$arr[17,3] = 'LO3924'
$arr[17,4] = 'Warsaw'
$arr[17,5] = '(WAW)'
$arr[17,6] = 'LOT '
$arr[17,7] = 'E190'
$arr[17,8] = '(SP-LMF)'
$arr[17,9] = '5:15 PM'
$arr[17,10] = $null
$arr[17,11] = '0 hours, 5 minutes'
$arr[17,12] = $null

$ws.Range("A188:M205").Value2 = $arr

Write-Output "Done. UsedRange: $($ws.UsedRange.Address())"